$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range('D2') '61.812.57'
Set-TextValue $ws.Range('E2') '  -0.22%  '

Set-TextValue $ws.Range('D3') '3.409.25'
Set-TextValue $ws.Range('E3') '  -0.19%  '

Set-TextValue $ws.Range('D4') '1.00'
Set-TextValue $ws.Range('E4') '  -0.09%  '

Set-TextValue $ws.Range('D5') '412.85'
Set-TextValue $ws.Range('E5') '  +0.98%  '

Set-TextValue $ws.Range('D6') '129.40'

Set-TextValue $ws.Range('D7') '0.622'
Set-TextValue $ws.Range('E7') '  -2.37%  '

Set-TextValue $ws.Range('E8') '  +0.06%  '

Set-TextValue $ws.Range('D9') '0.727'
Set-TextValue $ws.Range('E9') '  -0.80%  '

Set-TextValue $ws.Range('D10') '0.137'
Set-TextValue $ws.Range('E10') '  -3.96%  '

Set-TextValue $ws.Range('D11') '42.72'
Set-TextValue $ws.Range('E11') '  +0.80%  '

Set-TextValue $ws.Range('D12') '0.0000219'
Set-TextValue $ws.Range('E12') '  -0.67%  '

Set-TextValue $ws.Range('D13') '9.17'
Set-TextValue $ws.Range('E13') '  +2.18%  '

Set-TextValue $ws.Range('D14') '3.947.65'
Set-TextValue $ws.Range('E14') '  -0.31%  '

Set-TextValue $ws.Range('E15') '  -0.16%  '

Set-TextValue $ws.Range('D16') '20.46'
Set-TextValue $ws.Range('E16') '  -1.63%  '

Set-TextValue $ws.Range('B17') 'WrappedEther'
Set-TextValue $ws.Range('C17') 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws.Range('D17') '3.412.91'
Set-TextValue $ws.Range('E17') '  -0.11%  '

Set-TextValue $ws.Range('B18') 'Uniswap'
Set-TextValue $ws.Range('C18') 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue $ws.Range('D18') '12.70'
Set-TextValue $ws.Range('E18') '  +4.49%  '

Set-TextValue $ws.Range('E19') '  +0.76%  '

Set-TextValue $ws.Range('D20') '61.855.27'
Set-TextValue $ws.Range('E20') '  -0.13%  '

Set-TextValue $ws.Range('D21') '480.70'
Set-TextValue $ws.Range('E21') '  +9.69%  '

Set-TextValue $ws.Range('D22') '90.72'

Set-TextValue $ws.Range('D23') '3.28'
Set-TextValue $ws.Range('E23') '  +3.74%  '

Set-TextValue $ws.Range('D24') '13.10'
Set-TextValue $ws.Range('E24') '  +0.66%  '

Set-TextValue $ws.Range('D25') '3.32'
Set-TextValue $ws.Range('E25') '  +2.71%  '

Set-TextValue $ws.Range('D26') '9.80'
Set-TextValue $ws.Range('E26') '  +10.87%  '

Set-TextValue $ws.Range('D27') '33.07'
Set-TextValue $ws.Range('E27') '  -1.86%  '

Set-TextValue $ws.Range('E28') '  +0.42%  '

Set-TextValue $ws.Range('E29') '  +2.03%  '

Set-TextValue $ws.Range('D30') '11.87'
Set-TextValue $ws.Range('E30') '  -0.58%  '

Set-TextValue $ws.Range('E31') '  -1.63%  '

Set-TextValue $ws.Range('D32') '0.167'
Set-TextValue $ws.Range('E32') '  -1.72%  '

Set-TextValue $ws.Range('D33') '0.112'
Set-TextValue $ws.Range('E33') '  -3.03%  '

Set-TextValue $ws.Range('D34') '40.97'
Set-TextValue $ws.Range('E34') '  -3.54%  '

Set-TextValue $ws.Range('E35') '  +0.07%  '

Set-TextValue $ws.Range('D36') '57.82'
Set-TextValue $ws.Range('E36') '  +7.04%  '

Set-TextValue $ws.Range('D37') '0.0487'
Set-TextValue $ws.Range('E37') '  -2.43%  '

Set-TextValue $ws.Range('D38') '0.999'
Set-TextValue $ws.Range('E38') '  +0.08%  '

Set-TextValue $ws.Range('D39') '3.04'
Set-TextValue $ws.Range('E39') '  +4.47%  '

Set-TextValue $ws.Range('B40') 'TheGraph'
Set-TextValue $ws.Range('C40') 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextValue $ws.Range('D40') '0.327'
Set-TextValue $ws.Range('E40') '  +3.79%  '

Set-TextValue $ws.Range('B41') 'Monero'
Set-TextValue $ws.Range('C41') 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range('D41') '148.70'
Set-TextValue $ws.Range('E41') '  +5.43%  '

Set-TextValue $ws.Range('E42') '  -0.09%  '

Set-TextValue $ws.Range('D43') '3.33'
Set-TextValue $ws.Range('E43') '  -1.19%  '

Set-TextValue $ws.Range('D44') '2.08'
Set-TextValue $ws.Range('E44') '  +5.39%  '

Set-TextValue $ws.Range('B45') 'WEMIXToken'
Set-TextValue $ws.Range('C45') 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue $ws.Range('D45') '2.58'
Set-TextValue $ws.Range('E45') '  +6.76%  '

Set-TextValue $ws.Range('B46') 'NEARProtocol'
Set-TextValue $ws.Range('C46') 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Range('D46') '4.23'
Set-TextValue $ws.Range('E46') '  +3.14%  '

Set-TextValue $ws.Range('D47') '2.39'
Set-TextValue $ws.Range('E47') '  +21.17%  '

Set-TextValue $ws.Range('D48') '16.42'
Set-TextValue $ws.Range('E48') '  -0.75%  '

Set-TextValue $ws.Range('D49') '0.0₃0537'
Set-TextValue $ws.Range('E49') '  +20.51%  '

Set-TextValue $ws.Range('D50') '22.29'
Set-TextValue $ws.Range('E50') '  +0.45%  '

Set-TextValue $ws.Range('D51') '113.50'
Set-TextValue $ws.Range('E51') '  +11.09%  '
